# Add a new outstanding entry ("Anil Steel & Co.") as row 18 (Sr. No 9) on the
# "Purchase 22-23" sheet, pushing the existing row 18 entry ("Hatley
# Technologies", Sr. No 9) down to row 20 and renumbering it as Sr. No 10.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 18:19 - this shifts the old row 18 (and its
# relative formula/formatting) down to row 20, keeping row 19 blank just like
# the existing spacing pattern used throughout the sheet.
$ws.Rows("18:19").Insert()

# Give the newly freed row 18 the same look (fonts/borders/alignment/number
# formats) as the row it now sits above, then fill in its row height too.
$ws.Range("A20:F20").Copy()
$ws.Range("A18:F18").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows("18").RowHeight = 14.4

# Populate the new row 18 with the "Anil Steel & Co." entry.
$ws.Range("A18").Value = 9
$ws.Range("B18").Value = 45282
$ws.Range("C18").Value = 3262
$ws.Range("D18").Value = "Anil Steel & Co."
$ws.Range("E18").Value = 5629
$ws.Range("F18").Formula = "=E18"

# Renumber the pushed-down row as Sr. No 10.
$ws.Range("A20").Value = 10

# Match the selection left behind in the saved workbook.
$ws.Range("F18").Select()
